$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("LH-TC-REGISTERATION-Reviews")
$ws2 = $wb.Worksheets.Item("Version History")

# Row 6 of the Version History log: close out the reviewer-status entry.
$ws2.Range("B6").Value = "Ahmed Abuzaid"
$ws2.Range("C6").Value = "close reviewer status "
$ws2.Range("D6").Style = "Normal"

# Leave the selection where the author ended up when saving.
$ws1.Range("J9").Select()
$ws2.Range("C17").Select()
$ws2.Activate()
